# Auto-generated: apply cryptos-list price/volume refresh per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.977.02"
$ws.Range("E2").Value = "  +0.29%  "

$ws.Range("D3").Value = "3.210.97"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("E4").Value = "  +0.20%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "575.24"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -2.50%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "141.58"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -7.48%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.204.04"
$ws.Range("E8").Value = "  -0.47%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.524"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("E10").Value = "  -5.28%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.18"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.60%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.477"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.12%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.0000232"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -3.38%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "35.71"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -7.78%  "

$ws.Range("D15").Value = "3.723.46"
$ws.Range("E15").Value = "  -0.81%  "

$ws.Range("D16").Value = "66.914.81"
$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "3.215.64"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("E19").Value = "  -2.29%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "498.68"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -4.00%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.19"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.26%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.713"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -5.00%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "7.34"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.86%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "81.68"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -3.42%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "12.85"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.41%  "

$ws.Range("E26").Value = "  -0.09%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "3.05"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -4.42%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.03"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -5.47%  "

$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "27.67"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -3.68%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "7.56"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -3.70%  "

$ws.Range("E31").Value = "  +2.27%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "2.53"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.71%  "

$ws.Range("E33").Value = "  +0.01%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "510.87"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -4.88%  "

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "6.06"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -5.93%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "54.29"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "5.26"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -8.43%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0413"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -3.35%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.0810"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -4.58%  "

$ws.Range("E40").Value = "  -7.10%  "

$ws.Range("E41").Value = "  -5.21%  "

$ws.Range("D42").Value = "2.856.01"
$ws.Range("E42").Value = "  -1.60%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.51"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -11.98%  "

$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.09%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "0.249"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -3.33%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "121.73"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +1.68%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "24.73"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -6.12%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "2.01"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -6.24%  "

$ws.Range("D49").Value = "0.0₃0522"
$ws.Range("E49").Value = "  -11.22%  "

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.108"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -2.96%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.09"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -12.99%  "
